# edit.ps1 - apply the HAFIDH ASYI.docx edit described by the diff.
#
# Strategy: locate each affected paragraph via a unique text landmark using
# Find, then replace that paragraph's whole Range with precisely-crafted
# OOXML via Range.InsertXML. This keeps every untouched run/attribute byte
# for byte identical to the source, while giving full control over the
# handful of paragraphs that actually change.

$d = $word.ActiveDocument

function Get-ParagraphByLandmark($doc, [string]$landmark) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($landmark, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "landmark not found: $landmark"
    }
    return $rng.Paragraphs(1)
}

# ---------------------------------------------------------------------
# 1) Title block: "HAFIDH ASYI" paragraph -> empty para + Nama + NIM +
#    "BUG (LANJUTAN)" heading.
# ---------------------------------------------------------------------
$p = Get-ParagraphByLandmark $d "HAFIDH ASYI"
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Nama</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hafidh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Asyi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t>NIM</w:t></w:r><w:r><w:tab/><w:t>: 1957301050</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>BUG</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (LANJUTAN)</w:t></w:r></w:p>
</w:body>
</w:document>

'@
$p.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) Insert <w:lastRenderedPageBreak/> before "dengan" in
#    "...(low priority), bug dengan prioritas rendah merupakan...".
# ---------------------------------------------------------------------
$p = Get-ParagraphByLandmark $d "dengan prioritas rendah merupakan"
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00215494" w:rsidRDefault="008D0A30" w:rsidP="00F20875"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Priority </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>pengujian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>terbagi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>menjadi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>tiga</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>kelompok</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>yaitu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> low, medium </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> high priority. Bug yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>berdampak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>besar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>terhadap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>alur</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>bahkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>bisa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>merusak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>alur</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>adalah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> bug yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>sangat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>diprioritaskan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>perbaikannya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>atau</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>bisa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>disebut</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>sebagai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>tinggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> (high priority). </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>Setiap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>tinggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>harus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dilakukan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>perbaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>penyelesaian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>secepat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>mungkin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> agar proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>pengembangan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dapat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dilanjutkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>tidak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>terhambat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>Sedangkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>priotitas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>sedang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> (medium priority) yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>bisa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dibilang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>sebagai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> bug yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>tidak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>mengganggu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>alur</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>tidak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>harus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>langsung</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dilakukan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>perbaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve">, bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>sedang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dapat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>menunggu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>giliran</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>perbaikannya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>setelah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>tinggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>terselesaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>terlebih</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00875DA0"><w:t>dahulu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00875DA0"><w:t>.</w:t></w:r><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>Begitupun</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>halnya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>rendah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> (low priority), bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:lastRenderedPageBreak/><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>rendah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>merupakan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>kesalahan-kesalahan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>kecil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>tidak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>berdampak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>alur</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidR="00215494"><w:t>sama</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>sekali</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>Bahkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>priotitas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>rendah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>sering</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> kali </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>dibiarkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>begitu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>saja</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>tanpa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>adanya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00215494"><w:t>perbaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00215494"><w:t xml:space="preserve">. </w:t></w:r></w:p>
</w:body>
</w:document>

'@
$p.Range.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) Remove <w:lastRenderedPageBreak/> before "Proses" at the start of the
#    "Proses perbaikan bug dimulai..." paragraph.
# ---------------------------------------------------------------------
$p = Get-ParagraphByLandmark $d "Proses perbaikan bug dimulai"
$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="008D0A30" w:rsidRDefault="00215494" w:rsidP="00F20875"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perbaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dimulai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dari</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tinggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kemudian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dilanjutkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sedang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>diproses</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perbaikannya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>setelah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>semua</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tinggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>diselesaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>diakhiri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perbaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rendah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>diproses</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ketika</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>semua</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tinggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sedang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>selesai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>diperbaiki</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Perbaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rendah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tidak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>begitu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>disarankan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>karena</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hampir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>tidak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>berpengaruh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>Kesalahan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>kecil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>dalam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>pengkodean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidR="001D570F"><w:t>akan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>berdampak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>besar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>jalannya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>tersebut</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>bahkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>bisa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>merusak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>alur</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>sehingga</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>perbaikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>kesalahan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>tersebut</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>harus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>lebih</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>diprioritaskan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>walaupun</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>hanya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>kesalahan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00F20875"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00F20875"><w:t>kecil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>Berikut</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00AF1F41"><w:t>contoh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00AF1F41"><w:t xml:space="preserve"> bug </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00AF1F41"><w:t>berdasarkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00AF1F41"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>prioritas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D570F"><w:t>perbaikan</w:t></w:r><w:r w:rsidR="00AF1F41"><w:t>nya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D570F"><w:t>:</w:t></w:r></w:p>
</w:body>
</w:document>

'@
$p.Range.InsertXML($xml3)

# ---------------------------------------------------------------------
# 4) Blank paragraph before "Pengujian perangkat lunak tidak terlepas" ->
#    explicit page-break paragraph + "TEST CASE" heading (carries the
#    _GoBack bookmark and a lastRenderedPageBreak now that it starts a
#    new page).
# ---------------------------------------------------------------------
$p = Get-ParagraphByLandmark $d "Pengujian perangkat lunak tidak terlepas"
$blank = $p.Previous()
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>TEST CASE</w:t></w:r></w:p>
</w:body>
</w:document>

'@
$blank.Range.InsertXML($xml4)

# ---------------------------------------------------------------------
# 5) Remove <w:lastRenderedPageBreak/> before "akan" in "...apa yang akan
#    terjadi..." and remove the _GoBack bookmark after "...berhasil
#    login." (same paragraph, now relocated to the TEST CASE heading).
# ---------------------------------------------------------------------
$p = Get-ParagraphByLandmark $d "Cara penentuan/penulisan kasus uji"
$xml5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00AA0D49" w:rsidRDefault="00AA0D49" w:rsidP="00602024"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Cara </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>penentuan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>penulisan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kasus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>uji</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>baik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>didasarkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>beberapa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yaitu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>persyaratan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>awal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (preconditions), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>langkah-langkah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> proses (steps), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hasil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>diharapkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Persyaratan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>awal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>digunakan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>untuk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mendeskripsikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kesepakatan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>antara</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pihak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pengembang</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pelanggan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>terkait</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spesifikasi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>contohnya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>interaksi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pengguna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pada</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>halaman</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> login. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Langkah-langkah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> proses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>digunakan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>untuk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mendeskripsikan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lebih</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rinci</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>terkait</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>interaksi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>harus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dilakukan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oleh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pengguna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>perangkat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>lunak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>contohnya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>pengguna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>dapat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>melakukan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> login </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>dengan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> data-data yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>benar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve">/valid. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>Hasil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>diharapkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>digunakan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>untuk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>mengetahui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidR="006141D7"><w:t>apa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> yang </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>akan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>terjadi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>kemudian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>setelah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>langkah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>sebelumnya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>selesai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>dilakukan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>contohnya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>pengguna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>akan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>mendapatkan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> pop-up/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>pesan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>penyambutan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>setelah</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006141D7"><w:t>berhasil</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006141D7"><w:t xml:space="preserve"> login.</w:t></w:r></w:p>
</w:body>
</w:document>

'@
$p.Range.InsertXML($xml5)

Write-Output "done"
